# "add time stamp to excel sheet"
#
# The title-block label "Erstellt am" ("Created on") is replaced by
# "Datenaktualisierung:" ("Data update:") in cell D7 of both the
# Tabelle1 and Tabelle2 sheets (the "Datum"/"Date" value next to it in
# E7 is unchanged). The previously active sheet/selection (Tabelle3)
# is swapped for Tabelle1 becoming the active tab, with new selected
# cells recorded on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")
$ws3 = $wb.Worksheets.Item("Tabelle3")

# Rename the "Erstellt am" label to "Datenaktualisierung:" on both
# sheets that carry the title-block info table.
$ws1.Range("D7").Value = "Datenaktualisierung:"
$ws2.Range("D7").Value = "Datenaktualisierung:"

# Refresh the per-sheet selections to match the new view state.
$ws1.Range("G11").Select()
$ws2.Range("D8").Select()
$ws3.Range("C11").Select()

# Tabelle1 (first tab) becomes the active/selected sheet instead of
# Tabelle3.
$ws1.Activate()
